$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column cells: force text format to preserve exact formatting (e.g. trailing zeros)
$priceCells = @("D2","D3","D5","D6","D9","D14","D16","D17","D19","D20","D21","D24","D26","D29","D35","D36","D38","D41","D42","D43","D44","D45","D46","D47","D48")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply cell value updates
$ws.Range('D2').Value = '65.406.47'
$ws.Range('E2').Value = '  -0.80%  '
$ws.Range('D3').Value = '2.933.60'
$ws.Range('E3').Value = '  -2.72%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '568.33'
$ws.Range('E5').Value = '  -3.08%  '
$ws.Range('D6').Value = '158.61'
$ws.Range('E6').Value = '  +1.37%  '
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('E8').Value = '  -0.31%  '
$ws.Range('D9').Value = '2.928.70'
$ws.Range('E9').Value = '  -2.75%  '
$ws.Range('E10').Value = '  -4.03%  '
$ws.Range('E11').Value = '  -4.04%  '
$ws.Range('E12').Value = '  +0.30%  '
$ws.Range('E13').Value = '  -2.45%  '
$ws.Range('D14').Value = '34.20'
$ws.Range('E14').Value = '  -0.81%  '
$ws.Range('E15').Value = '  -0.81%  '
$ws.Range('D16').Value = '65.375.67'
$ws.Range('E16').Value = '  -0.80%  '
$ws.Range('D17').Value = '3.422.72'
$ws.Range('E17').Value = '  -2.65%  '
$ws.Range('E18').Value = '  -0.95%  '
$ws.Range('D19').Value = '2.934.60'
$ws.Range('E19').Value = '  -2.80%  '
$ws.Range('D20').Value = '15.11'
$ws.Range('E20').Value = '  +9.10%  '
$ws.Range('D21').Value = '444.72'
$ws.Range('E21').Value = '  -4.36%  '
$ws.Range('E22').Value = '  +0.07%  '
$ws.Range('E23').Value = '  -1.97%  '
$ws.Range('D24').Value = '82.09'
$ws.Range('E24').Value = '  -0.29%  '
$ws.Range('E25').Value = '  -1.73%  '
$ws.Range('D26').Value = '12.10'
$ws.Range('E26').Value = '  -3.44%  '
$ws.Range('E27').Value = '  -6.72%  '
$ws.Range('D29').Value = '8.07'
$ws.Range('E29').Value = '  +1.92%  '
$ws.Range('E30').Value = '  -0.61%  '
$ws.Range('E31').Value = '  -1.85%  '
$ws.Range('E32').Value = '  -5.62%  '
$ws.Range('E33').Value = '  -0.04%  '
$ws.Range('E34').Value = '  -1.03%  '
$ws.Range('D35').Value = '0.999'
$ws.Range('D36').Value = '0.971'
$ws.Range('E36').Value = '  -2.67%  '
$ws.Range('E37').Value = '  -1.87%  '
$ws.Range('D38').Value = '49.77'
$ws.Range('E38').Value = '  +1.06%  '
$ws.Range('E39').Value = '  -0.16%  '
$ws.Range('E40').Value = '  -10.28%  '
$ws.Range('B41').Value = 'TheGraph'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D41').Value = '0.299'
$ws.Range('E41').Value = '  -0.30%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').Value = '0.120'
$ws.Range('E42').Value = '  -2.53%  '
$ws.Range('D43').Value = '2.82'
$ws.Range('E43').Value = '  -7.13%  '
$ws.Range('D44').Value = '8.44'
$ws.Range('E44').Value = '  -0.85%  '
$ws.Range('D45').Value = '384.25'
$ws.Range('E45').Value = '  -2.28%  '
$ws.Range('D46').Value = '0.0350'
$ws.Range('E46').Value = '  -1.26%  '
$ws.Range('D47').Value = '2.699.18'
$ws.Range('E47').Value = '  -3.79%  '
$ws.Range('D48').Value = '133.17'
$ws.Range('E48').Value = '  -0.62%  '
$ws.Range('E49').Value = '  +0.01%  '
$ws.Range('E50').Value = '  +4.10%  '
$ws.Range('E51').Value = '  -1.63%  '
